$wb = $excel.ActiveWorkbook

# --- "survey" sheet: replace the old Ethiopia forms with the new "Demo" ones ---
$survey = $wb.Worksheets.Item("survey")

# Row 45: household-data form link now points at the Ethiopia_Demo_Household table
$survey.Range("B45").Value = "''?' + opendatakit.getHashString('../tables/Ethiopia_Demo_Household/forms/Ethiopia_household_data/',null)"

# Row 47/48: Section2 names form now lives under Ethiopia_Demo_Members (label text unchanged)
$survey.Range("A47").Value = "Ethiopia_Section2_names"
$survey.Range("B48").Value = "''?' + opendatakit.getHashString('../tables/Ethiopia_Demo_Members/forms/Ethiopia_Section2_names/',null)"

# Row 50/51: Section2 other_info form now lives under Ethiopia_Demo_Members (label text unchanged)
$survey.Range("A50").Value = "Ethiopia_Section2_other_info"
$survey.Range("B51").Value = "''?' + opendatakit.getHashString('../tables/Ethiopia_Demo_Members/forms/Ethiopia_Section2_other_info/',null)"

# Row 53/54: old Ethiopia_Section3 is replaced by the new Ethiopia_Section4 demo form
$survey.Range("A53").Value = "Ethiopia_Section4"
$survey.Range("B54").Value = "''?' + opendatakit.getHashString('../tables/Ethiopia_Demo_Members/forms/Ethiopia_Section4/',null)"

# Rows 53/54 shrink to match the new, shorter label/form rows
$survey.Rows(53).RowHeight = 16
$survey.Rows(54).RowHeight = 16

# The old Section 9.3 general/injuries forms (rows 56-61) are removed entirely
$survey.Rows("56:61").Delete()

# --- "choices" sheet: mirror the same form renames in the choice list ---
$choices = $wb.Worksheets.Item("choices")

$choices.Range("B15").Value = "Ethiopia_Section2_names"
$choices.Range("C15").Value = "Ethiopia_names"

$choices.Range("B16").Value = "Ethiopia_Section2_other_info"
$choices.Range("C16").Value = "Ethiopia_other_info"

$choices.Range("B17").Value = "Ethiopia_Section4"
$choices.Range("C17").Value = "Ethiopia Section 4"

# The old Section 9.3 general/injuries choice rows (18-20) are removed entirely
$choices.Rows("18:20").Delete()

# --- Window/selection state: focus moves from "choices" to "survey" ---
$survey.Activate()
$survey.Range("B45").Select()
